# DB schema v0.1 - "last updated and pending for mapping"
# Edits the "Sheet1" schema-reference worksheet:
#   - projects section header text: "project_details" -> "projects"
#   - phase status section header text: "project_status" -> "phase_status"
#   - role column of group_students now documents allowed values "0/1"
#   - widen column B to fit the new "0/1" note / longer type text
#   - move the active selection down to the phase_status header (A25)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the B63 value first so the new shared string "0/1" lands right after
# the existing "s_deptid" entry, matching the author's save order.
$ws.Range("B63").Value = "0/1"
$ws.Range("A12").Value = "projects"
$ws.Range("A25").Value = "phase_status"

# Column B needs to be noticeably wider now that it carries longer notes.
$ws.Columns.Item(2).ColumnWidth = 56.7

# Move/scroll the view's active selection to the phase_status header.
$ws.Range("A25").Select()
